$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of (cell address, new text value) pairs taken from the diff.
$updates = @(
    @("D2", "60.682.48"),
    @("E2", "  +4.94%  "),
    @("D3", "2.352.69"),
    @("E3", "  +2.81%  "),
    @("E4", "  +0.23%  "),
    @("D5", "546.88"),
    @("E5", "  +2.58%  "),
    @("D6", "132.36"),
    @("E6", "  +1.52%  "),
    @("E7", "  +0.14%  "),
    @("D8", "0.589"),
    @("E8", "  +1.71%  "),
    @("D9", "2.349.95"),
    @("E9", "  +2.73%  "),
    @("D10", "0.101"),
    @("E10", "  +2.00%  "),
    @("D11", "5.50"),
    @("E11", "  +1.60%  "),
    @("E12", "  +1.18%  "),
    @("D13", "0.335"),
    @("E13", "  +1.73%  "),
    @("D14", "23.97"),
    @("E14", "  +2.56%  "),
    @("D15", "2.771.77"),
    @("E15", "  +2.81%  "),
    @("D16", "60.678.55"),
    @("E16", "  +5.05%  "),
    @("D17", "0.0000133"),
    @("E17", "  +1.95%  "),
    @("D18", "2.355.48"),
    @("E18", "  +2.51%  "),
    @("D19", "10.71"),
    @("E19", "  +2.11%  "),
    @("D20", "4.18"),
    @("E20", "  -0.79%  "),
    @("D21", "6.88"),
    @("E21", "  +8.34%  "),
    @("D22", "315.29"),
    @("E22", "  +0.93%  "),
    @("D23", "0.997"),
    @("E23", "  -0.24%  "),
    @("D24", "63.41"),
    @("E24", "  +1.63%  "),
    @("D25", "0.172"),
    @("E25", "  +4.37%  "),
    @("D26", "1.00"),
    @("E26", "  +0.18%  "),
    @("D27", "7.96"),
    @("E27", "  -0.46%  "),
    @("E28", "  +6.49%  "),
    @("D29", "1.75"),
    @("E29", "  +3.12%  "),
    @("D30", "172.34"),
    @("E30", "  +0.83%  "),
    @("E31", "  +10.82%  "),
    @("D32", "0.0₃0731"),
    @("E32", "  +2.25%  "),
    @("D33", "5.90"),
    @("E33", "  +2.87%  "),
    @("E34", "  +15.90%  "),
    @("D35", "0.383"),
    @("E35", "  +1.05%  "),
    @("D36", "18.07"),
    @("E36", "  +2.20%  "),
    @("E37", "  +0.02%  "),
    @("E38", "  +0.18%  "),
    @("D39", "4.16"),
    @("E39", "  +7.23%  "),
    @("D40", "314.61"),
    @("E40", "  +9.49%  "),
    @("D41", "38.21"),
    @("E41", "  +0.32%  "),
    @("D42", "1.53"),
    @("E42", "  +3.54%  "),
    @("D43", "142.53"),
    @("E43", "  +0.94%  "),
    @("E44", "  +2.56%  "),
    @("D45", "0.0955"),
    @("E45", "  +1.00%  "),
    @("D46", "19.34"),
    @("E46", "  +7.46%  "),
    @("D47", "0.0499"),
    @("E47", "  +0.72%  "),
    @("D48", "0.562"),
    @("E48", "  +1.68%  "),
    @("D49", "0.0214"),
    @("E49", "  +2.21%  "),
    @("B50", "WhiteBITCoin"),
    @("C50", "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"),
    @("D50", "11.05"),
    @("E50", "  +1.01%  "),
    @("B51", "BabyDogeCoin"),
    @("C51", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"),
    @("D51", "0.0₆0211"),
    @("E51", "  +5.76%  ")
)

foreach ($pair in $updates) {
    $addr = $pair[0]
    $text = $pair[1]
    $cell = $ws.Range($addr)
    # Force the value to be written as text (matching the original inline-string cells)
    # instead of letting Excel auto-convert numeric-looking strings (e.g. "1.00") into numbers,
    # then restore the cell's original style so no formatting changes are introduced.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}
